$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 0.5888815306416006
$ws.Range("D2").Value = 0.03440194418249831
$ws.Range("E2").Value = 0.2072174822950821
$ws.Range("F2").Value = 0.9930695111709724
$ws.Range("G2").Value = 0.8482818272803314
$ws.Range("H2").Value = 0.8918630600367692
$ws.Range("I2").Value = 1.111877892119693
$ws.Range("K2").Value = 0.4984873897189175
$ws.Range("L2").Value = 0.181241298967322
$ws.Range("M2").Value = 0.1675686564766714
$ws.Range("N2").Value = 2.052689442208054

$ws.Range("B3").Value = 0.5706823736761351
$ws.Range("D3").Value = 0.03362665722089631
$ws.Range("E3").Value = 0.2070506507462229
$ws.Range("F3").Value = 0.9844754281123755
$ws.Range("G3").Value = 0.840836739305459
$ws.Range("H3").Value = 0.8930181420368513
$ws.Range("I3").Value = 1.119832566816203
$ws.Range("K3").Value = 0.4350155381087291
$ws.Range("L3").Value = 0.1743772246329343
$ws.Range("M3").Value = 0.161949245535272
$ws.Range("N3").Value = 2.072694104772502

$ws.Range("B4").Value = 0.5597590413339333
$ws.Range("D4").Value = 0.03314464114088267
$ws.Range("E4").Value = 0.2069566925088493
$ws.Range("F4").Value = 0.9798059982532976
$ws.Range("G4").Value = 0.8368245189965364
$ws.Range("H4").Value = 0.8941531576961381
$ws.Range("I4").Value = 1.125155115154378
$ws.Range("K4").Value = 0.3960058110855584
$ws.Range("L4").Value = 0.1702594723442559
$ws.Range("M4").Value = 0.1585776322827002
$ws.Range("N4").Value = 2.085608232925924

$ws.Range("B5").Value = 0.5553712386222855
$ws.Range("D5").Value = 0.03294672005085886
$ws.Range("E5").Value = 0.2069205369147475
$ws.Range("F5").Value = 0.9780558176451279
$ws.Range("G5").Value = 0.835329920727645
$ws.Range("H5").Value = 0.8947227818290457
$ws.Range("I5").Value = 1.127434378971664
$ws.Range("K5").Value = 0.38009981386449
$ws.Range("L5").Value = 0.1686058590775588
$ws.Range("M5").Value = 0.1572235676640119
$ws.Range("N5").Value = 2.091029553898888

$ws.Range("B6").Value = 0.5546464976476813
$ws.Range("D6").Value = 0.03291376536241586
$ws.Range("E6").Value = 0.2069146621928812
$ws.Range("F6").Value = 0.97777441890738
$ws.Range("G6").Value = 0.8350902194076895
$ws.Range("H6").Value = 0.8948238366546803
$ws.Range("I6").Value = 1.127819511434367
$ws.Range("K6").Value = 0.3774580866938777
$ws.Range("L6").Value = 0.1683327535476167
$ws.Range("M6").Value = 0.1569999307838401
$ws.Range("N6").Value = 2.091939345565564

$ws.Range("B7").Value = 0.5596996080447241
$ws.Range("D7").Value = 0.03314197795071649
$ws.Range("E7").Value = 0.2069561962617588
$ws.Range("F7").Value = 0.97978177673766
$ws.Range("G7").Value = 0.8368037940235098
$ws.Range("H7").Value = 0.8941604061841559
$ws.Range("I7").Value = 1.125185407501057
$ws.Range("K7").Value = 0.3957913341941151
$ws.Range("L7").Value = 0.1702370722473745
$ws.Range("M7").Value = 0.1585592902082524
$ws.Range("N7").Value = 2.08568070416276

$ws.Range("B8").Value = 0.582554607456899
$ws.Range("D8").Value = 0.03413587343716884
$ws.Range("E8").Value = 0.2071582001418841
$ws.Range("F8").Value = 0.9899801583003693
$ws.Range("G8").Value = 0.8455986110574969
$ws.Range("H8").Value = 0.892172973840232
$ws.Range("I8").Value = 1.114529719085237
$ws.Range("K8").Value = 0.4766102780446033
$ws.Range("L8").Value = 0.1788545044387746
$ws.Range("M8").Value = 0.165614800533497
$ws.Range("N8").Value = 2.059456044823769

$ws.Range("B9").Value = 0.6293492025062903
$ws.Range("D9").Value = 0.03603704987360601
$ws.Range("E9").Value = 0.2076215747410748
$ws.Range("F9").Value = 1.014804449880188
$ws.Range("G9").Value = 0.8672914390491968
$ws.Range("H9").Value = 0.8916538184188028
$ws.Range("I9").Value = 1.0971096446514
$ws.Range("K9").Value = 0.6347920171645001
$ws.Range("L9").Value = 0.1965199563278048
$ws.Range("M9").Value = 0.1800718485361301
$ws.Range("N9").Value = 2.013037285831651

$ws.Range("B10").Value = 0.66491676470946
$ws.Range("D10").Value = 0.03740432895317269
$ws.Range("E10").Value = 0.2080030670291322
$ws.Range("F10").Value = 1.035996537811627
$ws.Range("G10").Value = 0.885956767920959
$ws.Range("H10").Value = 0.8933328797678257
$ws.Range("I10").Value = 1.086426950043471
$ws.Range("K10").Value = 0.7508283742634774
$ws.Range("L10").Value = 0.2099657262155432
$ws.Range("M10").Value = 0.1910687859140339
$ws.Range("N10").Value = 1.981984513255263

$ws.Range("B11").Value = 0.6813518792844206
$ws.Range("D11").Value = 0.03801986571483695
$ws.Range("E11").Value = 0.2081855521918152
$ws.Range("F11").Value = 1.046281556856513
$ws.Range("G11").Value = 0.8950442564187284
$ws.Range("H11").Value = 0.8945444812609935
$ws.Range("I11").Value = 1.082025937063889
$ws.Range("K11").Value = 0.8035791507196564
$ws.Range("L11").Value = 0.2161839576054092
$ws.Range("M11").Value = 0.1961524452140395
$ws.Range("N11").Value = 1.968519951098166

$ws.Range("B12").Value = 0.687611763268734
$ws.Range("D12").Value = 0.03825201846254345
$ws.Range("E12").Value = 0.208255941200755
$ws.Range("F12").Value = 1.050269078512969
$ws.Range("G12").Value = 0.8985714828886273
$ws.Range("H12").Value = 0.8950676799717172
$ws.Range("I12").Value = 1.080425288942919
$ws.Range("K12").Value = 0.8235493843778841
$ws.Range("L12").Value = 0.2185532263191732
$ws.Range("M12").Value = 0.1980890747143675
$ws.Range("N12").Value = 1.963516412061956

$ws.Range("B13").Value = 0.6862619794880516
$ws.Range("D13").Value = 0.03820206208919785
$ws.Range("E13").Value = 0.2082407244792446
$ws.Range("F13").Value = 1.049406164946376
$ws.Range("G13").Value = 0.8978080036983584
$ws.Range("H13").Value = 0.8949521359272268
$ws.Range("I13").Value = 1.08076708575279
$ws.Range("K13").Value = 0.8192486802592498
$ws.Range("L13").Value = 0.2180423153721165
$ws.Range("M13").Value = 0.197671474503835
$ws.Range("N13").Value = 1.964589778506124

$ws.Range("B14").Value = 0.6818661594787443
$ws.Range("D14").Value = 0.03803898391554839
$ws.Range("E14").Value = 0.2081913173674232
$ws.Range("F14").Value = 1.046607752015078
$ws.Range("G14").Value = 0.8953327188544336
$ws.Range("H14").Value = 0.8945862346286333
$ws.Range("I14").Value = 1.081892929791593
$ws.Range("K14").Value = 0.8052222228693608
$ws.Range("L14").Value = 0.2163785873341624
$ws.Range("M14").Value = 0.1963115418417658
$ws.Range("N14").Value = 1.968106398568541

$ws.Range("B15").Value = 0.6791783021646154
$ws.Range("D15").Value = 0.03793897139290436
$ws.Range("E15").Value = 0.2081612215408617
$ws.Range("F15").Value = 1.044905732941032
$ws.Range("G15").Value = 0.8938277403310764
$ws.Range("H15").Value = 0.8943704952114757
$ws.Range("I15").Value = 1.082591125741089
$ws.Range("K15").Value = 0.7966299069818774
$ws.Range("L15").Value = 0.2153614000825286
$ws.Range("M15").Value = 0.1954800459422188
$ws.Range("N15").Value = 1.970272829836944

$ws.Range("B16").Value = 0.6638477898390249
$ws.Range("D16").Value = 0.03736397155313398
$ws.Range("E16").Value = 0.2079913211276363
$ws.Range("F16").Value = 1.035337368753218
$ws.Range("G16").Value = 0.88537490102172
$ws.Range("H16").Value = 0.8932627093967511
$ws.Range("I16").Value = 1.08672379253623
$ws.Range("K16").Value = 0.747380258401364
$ws.Range("L16").Value = 0.2095613909595073
$ws.Range("M16").Value = 0.1907381798495464
$ws.Range("N16").Value = 1.982877771515476

$ws.Range("B17").Value = 0.6545080785314212
$ws.Range("D17").Value = 0.03700956923712084
$ws.Range("E17").Value = 0.2078893829702926
$ws.Range("F17").Value = 1.029632671897218
$ws.Range("G17").Value = 0.8803422951296085
$ws.Range("H17").Value = 0.892697807817143
$ws.Range("I17").Value = 1.089376489195622
$ws.Range("K17").Value = 0.7171580505790587
$ws.Range("L17").Value = 0.2060292723392223
$ws.Range("M17").Value = 0.1878498958253587
$ws.Range("N17").Value = 1.990779973745269

$ws.Range("B18").Value = 0.6491601755564318
$ws.Range("D18").Value = 0.03680512037411887
$ws.Range("E18").Value = 0.2078315925332523
$ws.Range("F18").Value = 1.026412150335617
$ws.Range("G18").Value = 0.8775038096107579
$ws.Range("H18").Value = 0.892415039236468
$ws.Range("I18").Value = 1.090945417111271
$ws.Range("K18").Value = 0.6997717863231117
$ws.Range("L18").Value = 0.2040072657841421
$ws.Range("M18").Value = 0.1861962707090292
$ws.Range("N18").Value = 1.995387370233065

$ws.Range("B19").Value = 0.6473536119231085
$ws.Range("D19").Value = 0.03673579374796532
$ws.Range("E19").Value = 0.2078121702510147
$ws.Range("F19").Value = 1.025332153684545
$ws.Range("G19").Value = 0.8765523822317505
$ws.Range("H19").Value = 0.8923265376019884
$ws.Range("I19").Value = 1.091484043913127
$ws.Range("K19").Value = 0.6938845483271336
$ws.Range("L19").Value = 0.2033242954647676
$ws.Range("M19").Value = 0.1856376969004465
$ws.Range("N19").Value = 1.99695804627461

$ws.Range("B20").Value = 0.6554998206626976
$ws.Range("D20").Value = 0.03704735876825538
$ws.Range("E20").Value = 0.20790014735101
$ws.Range("F20").Value = 1.030233666777434
$ws.Range("G20").Value = 0.8808722131963833
$ws.Range("H20").Value = 0.8927535804071454
$ws.Range("I20").Value = 1.089089637677436
$ws.Range("K20").Value = 0.7203755950172024
$ws.Range("L20").Value = 0.2064042819074388
$ws.Range("M20").Value = 0.1881565686855993
$ws.Range("N20").Value = 1.989932326736628

$ws.Range("B21").Value = 0.6831563377313898
$ws.Range("D21").Value = 0.03808690946165427
$ws.Range("E21").Value = 0.2082057945359566
$ws.Range("F21").Value = 1.047427194057974
$ws.Range("G21").Value = 0.8960574345449146
$ws.Range("H21").Value = 0.8946919611373687
$ws.Range("I21").Value = 1.081560453260892
$ws.Range("K21").Value = 0.8093422811003279
$ws.Range("L21").Value = 0.2168668699480065
$ws.Range("M21").Value = 0.1967106741737652
$ws.Range("N21").Value = 1.967070897736509

$ws.Range("B22").Value = 0.7014426327668275
$ws.Range("D22").Value = 0.03876084552530301
$ws.Range("E22").Value = 0.2084130476185839
$ws.Range("F22").Value = 1.059205166363881
$ws.Range("G22").Value = 0.9064831781416842
$ws.Range("H22").Value = 0.8963341286532085
$ws.Range("I22").Value = 1.077023916458081
$ws.Range("K22").Value = 0.8674558642873649
$ws.Range("L22").Value = 0.2237896104822141
$ws.Range("M22").Value = 0.2023685976690714
$ws.Range("N22").Value = 1.952684558846681

$ws.Range("B23").Value = 0.6916637175719416
$ws.Range("D23").Value = 0.03840165726866474
$ws.Range("E23").Value = 0.2083017469192683
$ws.Range("F23").Value = 1.052869502002125
$ws.Range("G23").Value = 0.9008728219052813
$ws.Range("H23").Value = 0.8954233305700541
$ws.Range("I23").Value = 1.079410003748428
$ws.Range("K23").Value = 0.8364425281874048
$ws.Range("L23").Value = 0.2200870719107542
$ws.Range("M23").Value = 0.1993427317874037
$ws.Range("N23").Value = 1.960312017103398

$ws.Range("B24").Value = 0.6550513866641836
$ws.Range("D24").Value = 0.03703027629515532
$ws.Range("E24").Value = 0.2078952782370145
$ws.Range("F24").Value = 1.029961772641613
$ws.Range("G24").Value = 0.880632466470999
$ws.Range("H24").Value = 0.8927282347848546
$ws.Range("I24").Value = 1.089219186632221
$ws.Range("K24").Value = 0.7189209777760368
$ws.Range("L24").Value = 0.2062347131163733
$ws.Range("M24").Value = 0.1880179004229916
$ws.Range("N24").Value = 1.990315347610149

$ws.Range("B25").Value = 0.6164802444030784
$ws.Range("D25").Value = 0.03552788983231636
$ws.Range("E25").Value = 0.2074890155258062
$ws.Range("F25").Value = 1.007570994692003
$ws.Range("G25").Value = 0.8609451914199866
$ws.Range("H25").Value = 0.8914324997993077
$ws.Range("I25").Value = 1.10145045616024
$ws.Range("K25").Value = 0.5920316323767736
$ws.Range("L25").Value = 0.1916589481671593
$ws.Range("M25").Value = 0.1760946766563869
$ws.Range("N25").Value = 2.025058872275507

